# Update cryptos list (price/volume refresh + two coin-pair re-orderings)
# Leading "'" on column D values forces Excel to keep them as literal text
# (matching the original inlineStr cells) instead of auto-coercing them to
# numbers, which would otherwise drop things like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.198.50"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "'3.326.27"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'582.29"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("D6").Value = "'185.45"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'3.319.70"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'47.14"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'672.08"
$ws.Range("E14").Value = "  +9.62%  "
$ws.Range("D15").Value = "'3.858.80"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").Value = "'66.239.10"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.118"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.90"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'3.327.93"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("D24").Value = "'103.17"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "'3.98"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'9.52"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "'32.10"
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'599.21"
$ws.Range("E32").Value = "  +6.28%  "
$ws.Range("D33").Value = "'3.91"
$ws.Range("E33").Value = "  -5.07%  "
$ws.Range("D34").Value = "'10.99"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "'3.834.24"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'56.07"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").Value = "'2.69"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.127"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = ("'0.0{0}0698" -f [char]0x2083)
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").Value = "'32.83"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("E43").Value = "  +5.35%  "
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "'0.0413"
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("E47").Value = "  -12.95%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  +1.66%  "
